$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.019040600956404
$ws.Range("D2").Value = 1.039291936968012
$ws.Range("E2").Value = 1.033566015885498
$ws.Range("F2").Value = 1.045947579204553
$ws.Range("I2").Value = 1.057425463139092
$ws.Range("J2").Value = 1.040666352005827
$ws.Range("K2").Value = 1.050229723685538
$ws.Range("L2").Value = 1.044576720126875
$ws.Range("M2").Value = 1.056801859002528
$ws.Range("N2").Value = 1.017387091933935
$ws.Range("O2").Value = 1.03
$ws.Range("P2").Value = 1.053527136381266
$ws.Range("R2").Value = 1.046587007682587

$ws.Range("C3").Value = 1.022222178506418
$ws.Range("D3").Value = 1.041429850419536
$ws.Range("E3").Value = 1.035990814823534
$ws.Range("F3").Value = 1.048277440255041
$ws.Range("I3").Value = 1.05823519648247
$ws.Range("J3").Value = 1.042124845806073
$ws.Range("K3").Value = 1.051558712347278
$ws.Range("L3").Value = 1.046183000322007
$ws.Range("M3").Value = 1.058327732357383
$ws.Range("N3").Value = 1.017877386266501
$ws.Range("O3").Value = 1.03
$ws.Range("P3").Value = 1.054734740553169
$ws.Range("R3").Value = 1.047524060316341

$ws.Range("C4").Value = 1.024249349250888
$ws.Range("D4").Value = 1.042795792174699
$ws.Range("E4").Value = 1.037541418662536
$ws.Range("F4").Value = 1.049767848287503
$ws.Range("I4").Value = 1.058742910150015
$ws.Range("J4").Value = 1.043052456639351
$ws.Range("K4").Value = 1.052403753612854
$ws.Range("L4").Value = 1.047206919352071
$ws.Range("M4").Value = 1.059300550134616
$ws.Range("N4").Value = 1.018189274792223
$ws.Range("O4").Value = 1.03
$ws.Range("P4").Value = 1.055504645389015
$ws.Range("R4").Value = 1.04812244937176

$ws.Range("C5").Value = 1.025096703885108
$ws.Range("D5").Value = 1.043369423654775
$ws.Range("E5").Value = 1.038191327101152
$ws.Range("F5").Value = 1.050392785375392
$ws.Range("I5").Value = 1.058954470179815
$ws.Range("J5").Value = 1.04344106798788
$ws.Range("K5").Value = 1.052758944418609
$ws.Range("L5").Value = 1.047636063927537
$ws.Range("M5").Value = 1.059708476117422
$ws.Range("N5").Value = 1.018320220605798
$ws.Range("O5").Value = 1.03
$ws.Range("P5").Value = 1.05582748484421
$ws.Range("R5").Value = 1.048380683844607

$ws.Range("C6").Value = 1.025241584580902
$ws.Range("D6").Value = 1.043469733168212
$ws.Range("E6").Value = 1.038303050198004
$ws.Range("F6").Value = 1.050500390166448
$ws.Range("I6").Value = 1.05899215199467
$ws.Range("J6").Value = 1.043509059835752
$ws.Range("K6").Value = 1.052822565857439
$ws.Range("L6").Value = 1.047710728566107
$ws.Range("M6").Value = 1.059779648740752
$ws.Range("N6").Value = 1.01834346072552
$ws.Range("O6").Value = 1.03
$ws.Range("P6").Value = 1.055883812036371
$ws.Range("R6").Value = 1.048434315410396

$ws.Range("C7").Value = 1.024268841788764
$ws.Range("D7").Value = 1.04281480604089
$ws.Range("E7").Value = 1.037557773656845
$ws.Range("F7").Value = 1.049784009768475
$ws.Range("I7").Value = 1.058752137596863
$ws.Range("J7").Value = 1.043065610368448
$ws.Range("K7").Value = 1.052419723098175
$ws.Range("L7").Value = 1.047220237536333
$ws.Range("M7").Value = 1.059313726101495
$ws.Range("N7").Value = 1.01819459990847
$ws.Range("O7").Value = 1.03
$ws.Range("P7").Value = 1.055515073072029
$ws.Range("R7").Value = 1.048153628514919

$ws.Range("C8").Value = 1.020132584021978
$ws.Range("D8").Value = 1.040032072236245
$ws.Range("E8").Value = 1.034398870007013
$ws.Range("F8").Value = 1.04674820176494
$ws.Range("I8").Value = 1.057710372393234
$ws.Range("J8").Value = 1.041172475985893
$ws.Range("K8").Value = 1.05069577054959
$ws.Range("L8").Value = 1.045132205988629
$ws.Range("M8").Value = 1.057330091774744
$ws.Range("N8").Value = 1.017558380134057
$ws.Range("O8").Value = 1.03
$ws.Range("P8").Value = 1.053945189707959
$ws.Range("R8").Value = 1.046939191354186

$ws.Range("C9").Value = 1.012580846428425
$ws.Range("D9").Value = 1.034973205566912
$ws.Range("E9").Value = 1.028674752980116
$ws.Range("F9").Value = 1.041250746430632
$ws.Range("I9").Value = 1.055738083035235
$ws.Range("J9").Value = 1.037697032787753
$ws.Range("K9").Value = 1.047523252206095
$ws.Range("L9").Value = 1.041319123614355
$ws.Range("M9").Value = 1.053708223386654
$ws.Range("N9").Value = 1.016389403720052
$ws.Range("O9").Value = 1.03
$ws.Range("P9").Value = 1.05107877045623
$ws.Range("R9").Value = 1.044692920354108

$ws.Range("C10").Value = 1.007407938061314
$ws.Range("D10").Value = 1.03154902384844
$ws.Range("E10").Value = 1.024809835914185
$ws.Range("F10").Value = 1.037577588768751
$ws.Range("I10").Value = 1.054365839531409
$ws.Range("J10").Value = 1.035330693932132
$ws.Range("K10").Value = 1.04537109414728
$ws.Range("L10").Value = 1.038745138082557
$ws.Range("M10").Value = 1.051299969975319
$ws.Range("N10").Value = 1.01559809170883
$ws.Range("O10").Value = 1.03
$ws.Range("P10").Value = 1.049223600709562
$ws.Range("R10").Value = 1.043187879222651

$ws.Range("C11").Value = 1.005450277969986
$ws.Range("D11").Value = 1.030397856467636
$ws.Range("E11").Value = 1.023542155603758
$ws.Range("F11").Value = 1.036659018425292
$ws.Range("I11").Value = 1.053986580279953
$ws.Range("J11").Value = 1.03459806372062
$ws.Range("K11").Value = 1.04476841321441
$ws.Range("L11").Value = 1.038033350727398
$ws.Range("M11").Value = 1.050921200817051
$ws.Range("N11").Value = 1.015387066311509
$ws.Range("O11").Value = 1.03
$ws.Range("P11").Value = 1.049354651876338
$ws.Range("R11").Value = 1.042794269667602

$ws.Range("C12").Value = 1.004828267672657
$ws.Range("D12").Value = 1.030084932997003
$ws.Range("E12").Value = 1.023217189264317
$ws.Range("F12").Value = 1.036561409882244
$ws.Range("I12").Value = 1.053922267946464
$ws.Range("J12").Value = 1.034429682937646
$ws.Range("K12").Value = 1.04465759207884
$ws.Range("L12").Value = 1.037912704153967
$ws.Range("M12").Value = 1.051020143918623
$ws.Range("N12").Value = 1.015357044848001
$ws.Range("O12").Value = 1.03
$ws.Range("P12").Value = 1.049755629322875
$ws.Range("R12").Value = 1.04271591706007

$ws.Range("C13").Value = 1.005203830995843
$ws.Range("D13").Value = 1.030411615466078
$ws.Range("E13").Value = 1.023607519570107
$ws.Range("F13").Value = 1.037108435082689
$ws.Range("I13").Value = 1.054109728349106
$ws.Range("J13").Value = 1.034697754710775
$ws.Range("K13").Value = 1.044936324678749
$ws.Range("L13").Value = 1.038253430125241
$ws.Range("M13").Value = 1.051515863788566
$ws.Range("N13").Value = 1.015470104304432
$ws.Range("O13").Value = 1.03
$ws.Range("P13").Value = 1.050421500377913
$ws.Range("R13").Value = 1.042910519636174

$ws.Range("C14").Value = 1.005919863143914
$ws.Range("D14").Value = 1.030931677755966
$ws.Range("E14").Value = 1.024205355751062
$ws.Range("F14").Value = 1.037784269372285
$ws.Range("I14").Value = 1.054355468227715
$ws.Range("J14").Value = 1.035083731171145
$ws.Range("K14").Value = 1.045309164262682
$ws.Range("L14").Value = 1.038701179164815
$ws.Range("M14").Value = 1.052043236243294
$ws.Range("N14").Value = 1.015614350215626
$ws.Range("O14").Value = 1.03
$ws.Range("P14").Value = 1.051010198456026
$ws.Range("R14").Value = 1.043175533956055

$ws.Range("C15").Value = 1.006315579545782
$ws.Range("D15").Value = 1.031204604089819
$ws.Range("E15").Value = 1.024513728267523
$ws.Range("F15").Value = 1.03810137084214
$ws.Range("I15").Value = 1.054474716456778
$ws.Range("J15").Value = 1.035278448896669
$ws.Range("K15").Value = 1.045492303931229
$ws.Range("L15").Value = 1.038918254149191
$ws.Range("M15").Value = 1.052270685028348
$ws.Range("N15").Value = 1.01568294266227
$ws.Range("O15").Value = 1.03
$ws.Range("P15").Value = 1.051227214660328
$ws.Range("R15").Value = 1.043310826354418

$ws.Range("C16").Value = 1.00843237134864
$ws.Range("D16").Value = 1.032596855626875
$ws.Range("E16").Value = 1.026078208002621
$ws.Range("F16").Value = 1.039576809620459
$ws.Range("I16").Value = 1.055038003948437
$ws.Range("J16").Value = 1.03623965132487
$ws.Range("K16").Value = 1.046366465906088
$ws.Range("L16").Value = 1.039956708680342
$ws.Range("M16").Value = 1.053231768825802
$ws.Range("N16").Value = 1.016000427768509
$ws.Range("O16").Value = 1.03
$ws.Range("P16").Value = 1.051948453006411
$ws.Range("R16").Value = 1.043932012787276

$ws.Range("C17").Value = 1.009683281961035
$ws.Range("D17").Value = 1.033393846012458
$ws.Range("E17").Value = 1.026966697343574
$ws.Range("F17").Value = 1.040351553524935
$ws.Range("I17").Value = 1.055338868346042
$ws.Range("J17").Value = 1.036774365210657
$ws.Range("K17").Value = 1.046840642800227
$ws.Range("L17").Value = 1.040517865318297
$ws.Range("M17").Value = 1.053687225064935
$ws.Range("N17").Value = 1.016168173636774
$ws.Range("O17").Value = 1.03
$ws.Range("P17").Value = 1.052180528155815
$ws.Range("R17").Value = 1.044269828940008

$ws.Range("C18").Value = 1.010296614989336
$ws.Range("D18").Value = 1.033733383658372
$ws.Range("E18").Value = 1.027335065912811
$ws.Range("F18").Value = 1.040558174506773
$ws.Range("I18").Value = 1.055430781794463
$ws.Range("J18").Value = 1.036976193931656
$ws.Range("K18").Value = 1.04699417633374
$ws.Range("L18").Value = 1.040698071789966
$ws.Range("M18").Value = 1.053711712247134
$ws.Range("N18").Value = 1.01621596910208
$ws.Range("O18").Value = 1.03
$ws.Range("P18").Value = 1.051964506152659
$ws.Range("R18").Value = 1.044366867607193

$ws.Range("C19").Value = 1.010337923227231
$ws.Range("D19").Value = 1.033668898487524
$ws.Range("E19").Value = 1.027235189315402
$ws.Range("F19").Value = 1.040253393102796
$ws.Range("I19").Value = 1.055340712102965
$ws.Range("J19").Value = 1.03688360679892
$ws.Range("K19").Value = 1.046869145667099
$ws.Range("L19").Value = 1.040537689000802
$ws.Range("M19").Value = 1.053350664646683
$ws.Range("N19").Value = 1.01616019021899
$ws.Range("O19").Value = 1.03
$ws.Range("P19").Value = 1.051357128601772
$ws.Range("R19").Value = 1.044284795791685

$ws.Range("C20").Value = 1.008767200536705
$ws.Range("D20").Value = 1.032462030230044
$ws.Range("E20").Value = 1.025828047669021
$ws.Range("F20").Value = 1.038547787512597
$ws.Range("I20").Value = 1.054740419024287
$ws.Range("J20").Value = 1.035964295029898
$ws.Range("K20").Value = 1.045957120074344
$ws.Range("L20").Value = 1.039431452680544
$ws.Range("M20").Value = 1.051945073142759
$ws.Range("N20").Value = 1.015812702772913
$ws.Range("O20").Value = 1.03
$ws.Range("P20").Value = 1.049723711118798
$ws.Range("R20").Value = 1.043643849111855

$ws.Range("C21").Value = 1.004830786404262
$ws.Range("D21").Value = 1.029835118394689
$ws.Range("E21").Value = 1.022861743666727
$ws.Range("F21").Value = 1.035673962784843
$ws.Range("I21").Value = 1.053651746149892
$ws.Range("J21").Value = 1.034127913754824
$ws.Range("K21").Value = 1.04427247436777
$ws.Range("L21").Value = 1.037422562665723
$ws.Range("M21").Value = 1.050009596938229
$ws.Range("N21").Value = 1.015192246955193
$ws.Range("O21").Value = 1.03
$ws.Range("P21").Value = 1.048151592014068
$ws.Range("R21").Value = 1.042455954460747

$ws.Range("C22").Value = 1.002325722744493
$ws.Range("D22").Value = 1.028171715536657
$ws.Range("E22").Value = 1.020991564188131
$ws.Range("F22").Value = 1.033880417602086
$ws.Range("I22").Value = 1.052956127473325
$ws.Range("J22").Value = 1.032964175671502
$ws.Range("K22").Value = 1.04320500811888
$ws.Range("L22").Value = 1.036158336937532
$ws.Range("M22").Value = 1.0488093055857
$ws.Range("N22").Value = 1.014800000944778
$ws.Range("O22").Value = 1.03
$ws.Range("P22").Value = 1.047201651221383
$ws.Range("R22").Value = 1.041687853099489

$ws.Range("C23").Value = 1.003650338357008
$ws.Range("D23").Value = 1.029045771109199
$ws.Range("E23").Value = 1.02197832318119
$ws.Range("F23").Value = 1.03482645618998
$ws.Range("I23").Value = 1.053321156631309
$ws.Range("J23").Value = 1.033575958841048
$ws.Range("K23").Value = 1.043762966114749
$ws.Range("L23").Value = 1.036823575798381
$ws.Range("M23").Value = 1.049440595712084
$ws.Range("N23").Value = 1.015005244770884
$ws.Range("O23").Value = 1.03
$ws.Range("P23").Value = 1.047701270271637
$ws.Range("R23").Value = 1.042072808555805

$ws.Range("C24").Value = 1.008788719351419
$ws.Range("D24").Value = 1.032455067434425
$ws.Range("E24").Value = 1.025823571593902
$ws.Range("F24").Value = 1.038514839321518
$ws.Range("I24").Value = 1.054726171768082
$ws.Range("J24").Value = 1.03595260223616
$ws.Range("K24").Value = 1.045935239123694
$ws.Range("L24").Value = 1.039411880958234
$ws.Range("M24").Value = 1.051897743898245
$ws.Range("N24").Value = 1.01580430658382
$ws.Range("O24").Value = 1.03
$ws.Range("P24").Value = 1.049645916588001
$ws.Range("R24").Value = 1.043601256195377

$ws.Range("C25").Value = 1.014576912237285
$ws.Range("D25").Value = 1.036315430443237
$ws.Range("E25").Value = 1.030184570327555
$ws.Range("F25").Value = 1.042700775875183
$ws.Range("I25").Value = 1.056273146317765
$ws.Range("J25").Value = 1.038623470183979
$ws.Range("K25").Value = 1.048374989329401
$ws.Range("L25").Value = 1.042331707301915
$ws.Range("M25").Value = 1.054670503632368
$ws.Range("N25").Value = 1.016702483696514
$ws.Range("O25").Value = 1.03
$ws.Range("P25").Value = 1.05184033939951
$ws.Range("R25").Value = 1.045323419382978
